# Update tissue field with most detailed organ vs organ_part
# Target worksheet: "Tier 1_obs" (holds the sample/library metadata table)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tier 1_obs")

$ws.Range("S6").Value = 'blood'
$ws.Range("AB6").Value = 'SRR13806027, SRR13806059'
$ws.Range("AD6").Value = 'Illumina HiSeq 2500'
$ws.Range("S7").Value = 'kidney'
$ws.Range("AB7").Value = 'SRR13806028, SRR13806060'
$ws.Range("AD7").Value = 'Illumina HiSeq 2500'
$ws.Range("S8").Value = 'kidney'
$ws.Range("AB8").Value = 'SRR13806063'
$ws.Range("AD8").Value = 'Illumina HiSeq 2500'
$ws.Range("S9").Value = 'kidney'
$ws.Range("AB9").Value = 'SRR13806030, SRR13806062'
$ws.Range("AD9").Value = 'Illumina HiSeq 2500'
$ws.Range("S10").Value = 'kidney'
$ws.Range("AB10").Value = 'SRR13806061'
$ws.Range("AD10").Value = 'Illumina HiSeq 2500'
$ws.Range("S11").Value = 'blood'
$ws.Range("AB11").Value = 'SRR13806064'
$ws.Range("AD11").Value = 'Illumina HiSeq 2500'
$ws.Range("S12").Value = 'kidney'
$ws.Range("AB12").Value = 'SRR13806065'
$ws.Range("AD12").Value = 'Illumina HiSeq 2500'
$ws.Range("S13").Value = 'kidney'
$ws.Range("AB13").Value = 'SRR13806069'
$ws.Range("AD13").Value = 'Illumina HiSeq 2500'
$ws.Range("S14").Value = 'kidney'
$ws.Range("AB14").Value = 'SRR13806067'
$ws.Range("AD14").Value = 'Illumina HiSeq 2500'
$ws.Range("S15").Value = 'kidney'
$ws.Range("AB15").Value = 'SRR13806066'
$ws.Range("AD15").Value = 'Illumina HiSeq 2500'
$ws.Range("S16").Value = 'blood'
$ws.Range("AB16").Value = 'SRR13806070'
$ws.Range("AD16").Value = 'Illumina HiSeq 2500'
$ws.Range("S17").Value = 'kidney'
$ws.Range("AB17").Value = 'SRR13806071'
$ws.Range("AD17").Value = 'Illumina HiSeq 2500'
$ws.Range("S18").Value = 'lymph node'
$ws.Range("AB18").Value = 'SRR13806072'
$ws.Range("AD18").Value = 'Illumina HiSeq 2500'
$ws.Range("S19").Value = 'kidney'
$ws.Range("AB19").Value = 'SRR13806075'
$ws.Range("AD19").Value = 'Illumina HiSeq 2500'
$ws.Range("S20").Value = 'kidney'
$ws.Range("AB20").Value = 'SRR13806074'
$ws.Range("AD20").Value = 'Illumina HiSeq 2500'
$ws.Range("S21").Value = 'kidney'
$ws.Range("AB21").Value = 'SRR13806073'
$ws.Range("AD21").Value = 'Illumina HiSeq 2500'
$ws.Range("S22").Value = 'blood'
$ws.Range("AB22").Value = 'SRR13806076'
$ws.Range("AD22").Value = 'Illumina HiSeq 2500'
$ws.Range("S23").Value = 'kidney'
$ws.Range("AB23").Value = 'SRR13806077'
$ws.Range("AD23").Value = 'Illumina HiSeq 2500'
$ws.Range("S24").Value = 'kidney'
$ws.Range("AB24").Value = 'SRR13806078'
$ws.Range("AD24").Value = 'Illumina HiSeq 2500'
$ws.Range("S25").Value = 'kidney'
$ws.Range("AB25").Value = 'SRR13806023'
$ws.Range("AD25").Value = 'Illumina HiSeq 2500'
$ws.Range("S26").Value = 'kidney'
$ws.Range("AB26").Value = 'SRR13806024'
$ws.Range("AD26").Value = 'Illumina HiSeq 2500'
$ws.Range("S27").Value = 'kidney'
$ws.Range("AB27").Value = 'SRR13806057'
$ws.Range("AD27").Value = 'Illumina HiSeq 2500'
$ws.Range("S28").Value = 'kidney'
$ws.Range("AB28").Value = 'SRR13806080'
$ws.Range("AD28").Value = 'Illumina HiSeq 2500'
$ws.Range("S29").Value = 'kidney'
$ws.Range("AB29").Value = 'SRR13806068'
$ws.Range("AD29").Value = 'Illumina HiSeq 2500'
$ws.Range("S30").Value = 'kidney'
$ws.Range("AB30").Value = 'SRR13806079'
$ws.Range("AD30").Value = 'Illumina HiSeq 2500'
$ws.Range("S31").Value = 'kidney'
$ws.Range("AB31").Value = 'SRR13806025'
$ws.Range("AD31").Value = 'Illumina HiSeq 2500'
$ws.Range("S32").Value = 'kidney'
$ws.Range("AB32").Value = 'SRR13806046'
$ws.Range("AD32").Value = 'Illumina HiSeq 2500'
$ws.Range("S33").Value = 'kidney'
$ws.Range("AB33").Value = 'SRR13806035'
$ws.Range("AD33").Value = 'Illumina HiSeq 2500'
$ws.Range("S34").Value = 'kidney'
$ws.Range("AB34").Value = 'SRR13806026'
$ws.Range("AD34").Value = 'Illumina HiSeq 2500'
